# Card18: convert columns A:L (rows 2-12) from numeric/empty cells to text
# cells, matching a pandas `.astype(str)` style export (missing values
# become the literal string "nan"). Column M ("event") is already present
# and stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card18")

# Target text values for columns A (1) through L (12), rows 2-12.
$data = @{
    2  = @("2", "0", "150", "33", "✔", "nan", "nan", "nan", "nan", "nan", "nan", "nan")
    3  = @("2", "151", "300", "nan", "✔", "nan", "nan", "  ", "nan", "nan", "nan", "1\12\2024")
    4  = @("2", "301", "450", "nan", "nan", "nan", "nan", "nan", "nan", "nan", "nan", "nan")
    5  = @("2", "451", "550", "nan", "nan", "✔", "✔", "nan", "nan", "nan", "nan", "11\3\2025")
    6  = @("2", "551", "700", "590", "nan", "nan", "nan", "✔", "✔", "nan", "nan", "29\4\2025")
    7  = @("2", "701", "850", "785", "nan", "✔", "nan", "nan", "nan", "nan", "nan", "20\8\2025")
    8  = @("2", "851", "1000", "883", "nan", "✔", "✔", "nan", "nan", "nan", "nan", "23\10\2025")
    9  = @("2", "1001", "1150", "nan", "nan", "nan", "nan", "nan", "nan", "nan", "nan", "nan")
    10 = @("2", "1151", "1300", "nan", "nan", "nan", "nan", "nan", "nan", "nan", "nan", "nan")
    11 = @("2", "1301", "1450", "nan", "nan", "nan", "nan", "nan", "nan", "nan", "nan", "nan")
    12 = @("2", "1451", "1500", "nan", "nan", "nan", "nan", "nan", "nan", "nan", "nan", "nan")
}

foreach ($r in 2..12) {
    $rowValues = $data[$r]
    for ($c = 1; $c -le 12; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        # Prefix with an apostrophe so Excel stores the value as literal
        # text instead of re-parsing numeric-looking strings (e.g. "2",
        # "150") back into numbers; ClearFormats drops the transient
        # quote-prefix style Excel applies so no stray style/format is
        # left behind on the cell.
        $cell.Value = "'" + $rowValues[$c - 1]
        $cell.ClearFormats()
    }
}
